$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 1.04489337822671
$ws.Range("C4").Value = 1.8131868131868101
$ws.Range("D4").Value = 1.5283505154639201
$ws.Range("B5").Value = 1.06904487917146
$ws.Range("C5").Value = 1.81938325991189
$ws.Range("D5").Value = 1.5416666666666701
$ws.Range("B6").Value = 1.83769633507853
$ws.Range("C6").Value = 1.9842931937172801
$ws.Range("D6").Value = 3.5241379310344798
$ws.Range("B7").Value = 1.1829652996845399
$ws.Range("C7").Value = 1.32547169811321
$ws.Range("D7").Value = 2.3734439834024901
$ws.Range("B8").Value = 1.0407876230661
$ws.Range("C8").Value = 1.2452554744525499
$ws.Range("D8").Value = 2.22900763358779
$ws.Range("B9").Value = 1.9939577039274901
$ws.Range("C9").Value = 2.46875
$ws.Range("D9").Value = 2.8443113772455102
$ws.Range("B10").Value = 0.32226720647773299
$ws.Range("C10").Value = 0.40614886731391597
$ws.Range("D10").Value = 1.87341772151899
$ws.Range("B11").Value = 0.27869986168741401
$ws.Range("C11").Value = 0.36480686695279002
$ws.Range("D11").Value = 1.87341772151899
$ws.Range("B12").Value = 0.34806629834254099
$ws.Range("C12").Value = 0.42547660311958402
$ws.Range("D12").Value = 1.87341772151899
$ws.Range("B13").Value = 0.35271317829457399
$ws.Range("C13").Value = 0.43416370106761598
$ws.Range("D13").Value = 1.87341772151899
$ws.Range("B14").Value = 0.32226720647773299
$ws.Range("C14").Value = 0.40614886731391597
$ws.Range("D14").Value = 1.87341772151899
$ws.Range("B15").Value = 0.53144266337854495
$ws.Range("C15").Value = 0.64546525323910497
$ws.Range("D15").Value = 1.71428571428571
$ws.Range("B16").Value = 0.53144266337854495
$ws.Range("C16").Value = 0.64546525323910497
$ws.Range("D16").Value = 1.71428571428571
$ws.Range("B17").Value = 0.60239361702127703
$ws.Range("C17").Value = 0.61619718309859195
$ws.Range("D17").Value = 1.4833836858006
$ws.Range("B18").Value = 0.60397350993377497
$ws.Range("C18").Value = 0.61592505854800905
$ws.Range("D18").Value = 1.4700598802395199
$ws.Range("B19").Value = 0.61272475795297399
$ws.Range("C19").Value = 0.64444444444444504
$ws.Range("D19").Value = 1.4787878787878801
$ws.Range("B20").Value = 0.54235294117647104
$ws.Range("C20").Value = 0.55835962145110396
$ws.Range("D20").Value = 1.2517482517482501
$ws.Range("B21").Value = 0.54427390791027197
$ws.Range("C21").Value = 0.56118143459915604
$ws.Range("D21").Value = 1.24708624708625
$ws.Range("B22").Value = 0.5625
$ws.Range("C22").Value = 0.58914728682170603
$ws.Range("D22").Value = 1.24942263279446
$ws.Range("B23").Value = 0.49432739059967601
$ws.Range("C23").Value = 0.52307692307692299
$ws.Range("D23").Value = 0.86149584487534603
$ws.Range("B24").Value = 0.47744945567651598
$ws.Range("C24").Value = 0.50454545454545496
$ws.Range("D24").Value = 0.83858267716535395
$ws.Range("B25").Value = 0.49833333333333302
$ws.Range("C25").Value = 0.53212121212121199
$ws.Range("D25").Value = 0.89019033674963399
$ws.Range("B26").Value = 0.61425576519916203
$ws.Range("C26").Value = 0.64909638554216897
$ws.Range("D26").Value = 1.1215932914046101
$ws.Range("B27").Value = 0.495207667731629
$ws.Range("C27").Value = 0.52470588235294102
$ws.Range("D27").Value = 0.84471218206158005
$ws.Range("B28").Value = 0.811594202898551
$ws.Range("C28").Value = 0.81545741324921095
$ws.Range("D28").Value = 2.1052631578947398
$ws.Range("B29").Value = 0.44981862152357899
$ws.Range("C29").Value = 0.47120921305182301
$ws.Range("D29").Value = 1.01876675603217
$ws.Range("B30").Value = 0.43322109988776702
$ws.Range("C30").Value = 0.46581586144029202
$ws.Range("D30").Value = 0.95971563981042696

$ws.Range("C6").Select()
